# B6-PowerPoint.pptx edit — Tue, Jun 16, 2020  1:05:44 PM
#
# 1) Re-style the three data tables (slides 14, 15, 16) from the custom
#    "Table_0" style ({F4CDD322-E8DC-4F45-A7AF-51536E0D2228}) to the
#    built-in "No Style, No Grid" table style
#    ({A58F5878-5EA5-4C03-9E9B-FAC01E448CB0}).
# 2) Re-colour the deck's theme palette from the "Integral" (Red Violet)
#    scheme over to the stock "Office" scheme.

$p = $ppt.ActivePresentation

# --- 1) Tables --------------------------------------------------------
$newStyleId = "{A58F5878-5EA5-4C03-9E9B-FAC01E448CB0}"
$tableSlideIndexes = @(14, 15, 16)
foreach ($idx in $tableSlideIndexes) {
    $slide = $p.Slides.Item($idx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2) Theme colours ---------------------------------------------------
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink -> stock Office palette
$officeColors = @(0x000000, 0xFFFFFF, 0x6A5444, 0xE6E6E7, 0xD59B5B, 0x317DED, 0xA5A5A5, 0x00C0FF, 0xC47244, 0x47AD70, 0xC16305, 0x724F95)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i - 1]
}
